$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

# Add new row 12 with the new change-track entry
$ws.Cells.Item(12, 1).Value = (Get-Date -Year 2014 -Month 3 -Day 29)
$ws.Cells.Item(12, 1).Style = $ws.Cells.Item(11, 1).Style
$ws.Cells.Item(12, 1).NumberFormat = $ws.Cells.Item(11, 1).NumberFormat

$ws.Cells.Item(12, 2).Value = "10"
$ws.Cells.Item(12, 2).Style = $ws.Cells.Item(11, 2).Style
$ws.Cells.Item(12, 2).NumberFormat = $ws.Cells.Item(11, 2).NumberFormat

$ws.Cells.Item(12, 3).Value = "SPA"
$ws.Cells.Item(12, 3).Style = $ws.Cells.Item(11, 3).Style
$ws.Cells.Item(12, 3).NumberFormat = $ws.Cells.Item(11, 3).NumberFormat

$ws.Cells.Item(12, 4).Value = "OS Dispatcher and OSEK functions"

$ws.Cells.Item(12, 5).Value = "In process"

# Update the active selection like in the diff
$ws.Range("D21").Select()
